# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" (positioned right before "总计"),
#   populated with the Q1-2022 fund-holding detail rows.
# - Update the "总计" (totals) worksheet: add a new top row summarizing
#   2022-Q1 and push the existing 2021-Q4 / 2021-Q3 rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the "2022-Q1" sheet by duplicating "总计" (so it inherits the
#    exact same header/row styling), inserted immediately before "总计".
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Copy($total)

# NOTE: worksheet object references are position-bound, so after the
# insert-by-copy above, the old "$total" variable now actually points at
# the freshly inserted copy. Re-resolve both sheets by name so the rest
# of the script talks to the correct objects.
$newSheet = $wb.Worksheets.Item("总计 (2)")
$total = $wb.Worksheets.Item("总计")
$newSheet.Name = "2022-Q1"

# Wipe the copied values but keep the inherited cell formatting.
$newSheet.Range("A1:D3").ClearContents()
# The new sheet only needs a header row + a single data row.
$newSheet.Rows.Item(3).Delete()

# Extend the header styling (copied from column D) across E:H.
$newSheet.Range("D1").Copy()
$newSheet.Range("E1:H1").PasteSpecial(-4122)

# --- Header row ---
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# --- Data row (fund codes / percentages are stored as text, matching
#     the other quarterly sheets) ---
$newSheet.Range("A2").Value = 0

$c = $newSheet.Range("B2")
$c.NumberFormat = "@"
$c.Value = "000049"
$c.Style = "Normal"

$newSheet.Range("C2").Value = "中银标普全球精选自然资源等权重指数(QDII)"

$c = $newSheet.Range("D2")
$c.NumberFormat = "@"
$c.Value = "0.27"
$c.Style = "Normal"

$c = $newSheet.Range("E2")
$c.NumberFormat = "@"
$c.Value = "89.72"
$c.Style = "Normal"

$c = $newSheet.Range("F2")
$c.NumberFormat = "@"
$c.Value = "1.26"
$c.Style = "Normal"

$c = $newSheet.Range("G2")
$c.NumberFormat = "@"
$c.Value = "0.0034"
$c.Style = "Normal"

$newSheet.Range("H2").Value = 4

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: push the existing two rows down by one and
#    insert the new 2022-Q1 summary row at the top of the data.
# ---------------------------------------------------------------------
$b2 = $total.Range("B2").Value()
$c2 = $total.Range("C2").Value()
$d2 = $total.Range("D2").Value()
$b3 = $total.Range("B3").Value()
$c3 = $total.Range("C3").Value()
$d3 = $total.Range("D3").Value()

# Give the new row 4 (index column) the same style as the existing index
# column cells before writing into it.
$total.Range("A2").Copy()
$total.Range("A4").PasteSpecial(-4122)

# 2021-Q3 row moves from row 3 to row 4
$total.Range("A4").Value = 2
$total.Range("B4").Value = $b3
$total.Range("C4").Value = $c3
$total.Range("D4").Value = $d3

# 2021-Q4 row moves from row 2 to row 3
$total.Range("A3").Value = 1
$total.Range("B3").Value = $b2
$total.Range("C3").Value = $c2
$total.Range("D3").Value = $d2

# New 2022-Q1 summary row becomes row 2
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0
